# diary.xlsx - weekly update: 7 marras (7 Nov) entry gets a second block of
# hours, a "learning quality" note, a code/META note, and the book-progress
# note is extended with the chapter reference. Tunnit (G26) bumps 2 -> 2.5,
# which ripples into the Kertyma total (H3) automatically via SUM(G3:G60).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Write the new / changed cell contents for row 26 ("7 marras").
# Order matters for how new entries land in the shared-string table, so we
# write B (extended time range) first, then the two brand-new narrative
# cells (D, F), and finally the C cell whose text was tweaked.
$ws.Range("B26").Value = "12.45-14.45, 15.45-16.15"
$ws.Range("D26").Value = "Ihan Jees, suht simppelisti rakentuu sen partikkelijärjestelmän päälle. Katsotaan vielä ensi viikolla mitä jäi käteen teoriatankkauksesta ja täydennetään aukkoja tekemällä demo."
$ws.Range("F26").Value = "Partikkeleilleko voisi siis saada vesisimun aikaan? Voisi olla hyvä niitä kypsytellä vielä kun kirjasta alkaa olla enemmän takana kuin edessä."
$ws.Range("C26").Value = "Kovat kappaleet, s.194-213 luku 10."

# Tunnit (hours) for the day: 2 -> 2.5
$ws.Range("G26").Value = 2.5

# Match formatting used by sibling rows: B column keeps the time number
# format with wrap text (like B23/B25 etc.), the new narrative cells (D, F)
# use the same wrap-text "General" style as the rest of the table's prose
# columns.
$ws.Range("B26").NumberFormat = "h:mm"
$ws.Range("B26").WrapText = $true
$ws.Range("D26").WrapText = $true
$ws.Range("F26").WrapText = $true

# The row now holds much more text, so it grows taller to fit (matches the
# other wrapped, multi-column rows in the sheet).
$ws.Rows.Item(26).RowHeight = 72.5

# Scroll the view down to the bottom of the log and leave E26 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 24
$win.ScrollColumn = 1
[void]$ws.Range("E26").Select()
